$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H49").Value = 2414.5454
$ws.Range("I49").Value = 650
$ws.Range("J49").Value = 3422.8572
$ws.Range("K49").Value = 1950
$ws.Range("L49").Value = 10268.5716
$ws.Range("M49").Value = -1814
$ws.Range("N49").Value = -10540.5716
$ws.Range("H98").Value = 1388.9117
$ws.Range("I98").Value = 1150.3928
$ws.Range("K98").Value = 1150.3928
$ws.Range("M98").Value = 347.6071999999999
$ws.Range("H111").Value = 4104.2856
$ws.Range("I111").Value = 3298.5
$ws.Range("J111").Value = 5178.6665
$ws.Range("K111").Value = 9895.5
$ws.Range("L111").Value = 15535.9995
$ws.Range("M111").Value = -6828.5
$ws.Range("N111").Value = -21669.9995
$ws.Range("H113").Value = 2126.1428
$ws.Range("I113").Value = 2319.375
$ws.Range("J113").Value = 2007.2307
$ws.Range("K113").Value = 2319.375
$ws.Range("L113").Value = 2007.2307
$ws.Range("M113").Value = 934.625
$ws.Range("N113").Value = -8515.2307
$ws.Range("H122").Value = 1388.9117
$ws.Range("I122").Value = 1150.3928
$ws.Range("K122").Value = 3451.1784
$ws.Range("M122").Value = -1001.1784
$ws.Range("H140").Value = 47241.2
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 47241.2
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 47241.2
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -57601.2

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1953.08
$ws.Range("I132").Value = 1680.6757
$ws.Range("J132").Value = 2728.3845
$ws.Range("K132").Value = 5042.0271
$ws.Range("L132").Value = 8185.1535
$ws.Range("M132").Value = -2512.0271
$ws.Range("N132").Value = -13245.1535

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 3373176.2
$ws.Range("J9").Value = 3373176.2
$ws.Range("L9").Value = 3373176.2
$ws.Range("N9").Value = -3373512.2
$ws.Range("H81").Value = 18037.375
$ws.Range("J81").Value = 18037.375
$ws.Range("L81").Value = 18037.375
$ws.Range("N81").Value = -20159.375
$ws.Range("H84").Value = 18037.375
$ws.Range("J84").Value = 18037.375
$ws.Range("L84").Value = 54112.125
$ws.Range("N84").Value = -64720.125

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 2332
$ws.Range("I59").Value = 1998
$ws.Range("J59").Value = 3000
$ws.Range("K59").Value = 5994
$ws.Range("L59").Value = 9000
$ws.Range("M59").Value = -5454
$ws.Range("N59").Value = -10080
$ws.Range("H68").Value = 987.75
$ws.Range("I68").Value = 734
$ws.Range("K68").Value = 2202
$ws.Range("M68").Value = -1391
$ws.Range("H71").Value = 987.75
$ws.Range("I71").Value = 734
$ws.Range("K71").Value = 6606
$ws.Range("M71").Value = -2550
$ws.Range("H104").Value = 4089.4285
$ws.Range("I104").Value = 763
$ws.Range("J104").Value = 5420
$ws.Range("K104").Value = 2289
$ws.Range("L104").Value = 16260
$ws.Range("M104").Value = 332
$ws.Range("N104").Value = -21502
$ws.Range("H105").Value = 12311.111
$ws.Range("J105").Value = 12311.111
$ws.Range("L105").Value = 36933.333
$ws.Range("N105").Value = -42175.333
$ws.Range("H106").Value = 19000
$ws.Range("J106").Value = 19000
$ws.Range("L106").Value = 57000
$ws.Range("N106").Value = -58892
$ws.Range("H113").Value = 562375
$ws.Range("I113").Value = 517.625
$ws.Range("J113").Value = 936946.5600000001
$ws.Range("K113").Value = 1552.875
$ws.Range("L113").Value = 2810839.68
$ws.Range("M113").Value = 617.125
$ws.Range("N113").Value = -2815179.68
$ws.Range("H139").Value = 19790.643
$ws.Range("I139").Value = 1434.2433
$ws.Range("J139").Value = 55537.316
$ws.Range("K139").Value = 4302.7299
$ws.Range("L139").Value = 166611.948
$ws.Range("M139").Value = 837.2700999999997
$ws.Range("N139").Value = -176891.948
$ws.Range("H140").Value = 25785.047
$ws.Range("I140").Value = 59940.53
$ws.Range("J140").Value = 3452.6155
$ws.Range("K140").Value = 179821.59
$ws.Range("L140").Value = 10357.8465
$ws.Range("M140").Value = -174641.59
$ws.Range("N140").Value = -20717.8465

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 598.74194
$ws.Range("I107").Value = 464.61905
$ws.Range("J107").Value = 880.4
$ws.Range("K107").Value = 464.61905
$ws.Range("L107").Value = 880.4
$ws.Range("M107").Value = 1455.38095
$ws.Range("N107").Value = -4720.4
$ws.Range("H130").Value = 32828.57
$ws.Range("J130").Value = 32828.57
$ws.Range("L130").Value = 32828.57
$ws.Range("N130").Value = -42868.57
$ws.Range("H132").Value = 1776.3334
$ws.Range("I132").Value = 1535.5714
$ws.Range("J132").Value = 2899.889
$ws.Range("K132").Value = 4606.7142
$ws.Range("L132").Value = 8699.667000000001
$ws.Range("M132").Value = -2076.7142
$ws.Range("N132").Value = -13759.667
$ws.Range("H133").Value = 25374.637
$ws.Range("J133").Value = 25374.637
$ws.Range("L133").Value = 25374.637
$ws.Range("N133").Value = -35494.637
$ws.Range("H135").Value = 62673.332
$ws.Range("J135").Value = 62673.332
$ws.Range("L135").Value = 62673.332
$ws.Range("N135").Value = -72813.33199999999
$ws.Range("H137").Value = 50780
$ws.Range("J137").Value = 50780
$ws.Range("L137").Value = 50780
$ws.Range("N137").Value = -60980
$ws.Range("H138").Value = 82666.664
$ws.Range("J138").Value = 82666.664
$ws.Range("L138").Value = 82666.664
$ws.Range("N138").Value = -92946.664
$ws.Range("H139").Value = 45092
$ws.Range("J139").Value = 45092
$ws.Range("L139").Value = 45092
$ws.Range("N139").Value = -55372
$ws.Range("H140").Value = 41052.875
$ws.Range("J140").Value = 41052.875
$ws.Range("L140").Value = 41052.875
$ws.Range("N140").Value = -51412.875
$ws.Range("H141").Value = 50676.145
$ws.Range("J141").Value = 50676.145
$ws.Range("L141").Value = 50676.145
$ws.Range("N141").Value = -61036.145

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4023.3635
$ws.Range("I132").Value = 4977.1113
$ws.Range("J132").Value = 2878.8667
$ws.Range("K132").Value = 14931.3339
$ws.Range("L132").Value = 8636.6001
$ws.Range("M132").Value = -12401.3339
$ws.Range("N132").Value = -13696.6001

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H132").Value = 2494.3428
$ws.Range("I132").Value = 3213.55
$ws.Range("J132").Value = 1535.4
$ws.Range("K132").Value = 9640.650000000001
$ws.Range("L132").Value = 4606.200000000001
$ws.Range("M132").Value = -7110.650000000001
$ws.Range("N132").Value = -9666.200000000001
